$wb = $excel.ActiveWorkbook

# --- Add a brand-new "Player Info" worksheet in front of "ODI Batting" ---
$battingSheetForAnchor = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetForAnchor)
$playerInfo.Name = "Player Info"

# Re-fetch the other sheet references *after* inserting the new sheet, since
# handles captured before an insert can point at the wrong sheet afterwards.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the header styling used on the other sheets (bold, centered/top, thin border).
# Style cell-by-cell (border/alignment before font) so the engine reuses a single
# cached style combo for every header cell instead of minting one per cell.
foreach ($addr in @("A1", "B1", "C1", "D1")) {
    $cell = $playerInfo.Range($addr)
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Font.Bold = $true
}

# Data row (leading "'" forces text storage for the numeric-looking id, matching
# how the other "numeric" text columns in this workbook are stored)
$playerInfo.Range("A2").Value = "'4670"
$playerInfo.Range("B2").Value = "Rahul Desraj Chahar"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

$playerInfo.Range("A1").Select() | Out-Null

# --- ODI Batting: MATCH_CARD_LINK -> MATCH_CODE, url -> bare match code ---
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").Value = "'4485"

# --- ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE, url -> bare match code ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").Value = "'4485"

Write-Host "Player Info sheet added; ODI Batting/Bowling MATCH_CODE columns updated."
